# The deck's theme was swapped: the slide master (and therefore every
# slide) switches from the custom "TU Braunschweig" color scheme back to
# the generic Office "Default" color scheme (the palette that used to sit
# unused in the deck's secondary theme part).
#
# PowerPoint's object model exposes the live, rendered color scheme via
# SlideMaster.ColorScheme - each of the twelve theme slots (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) is reachable as Colors(1..12) and each
# one carries a settable .RGB (VBA-style 0xBBGGRR integer).

function ColorValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$cs.Colors(1).RGB  = ColorValue 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = ColorValue 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = ColorValue 0x15 0x81 0x58   # dk2
$cs.Colors(4).RGB  = ColorValue 0xF3 0xF3 0xF3   # lt2
$cs.Colors(5).RGB  = ColorValue 0x05 0x8D 0xC7   # accent1
$cs.Colors(6).RGB  = ColorValue 0x50 0xB4 0x32   # accent2
$cs.Colors(7).RGB  = ColorValue 0xED 0x56 0x1B   # accent3
$cs.Colors(8).RGB  = ColorValue 0xED 0xEF 0x00   # accent4
$cs.Colors(9).RGB  = ColorValue 0x24 0xCB 0xE5   # accent5
$cs.Colors(10).RGB = ColorValue 0x64 0xE5 0x72   # accent6
$cs.Colors(11).RGB = ColorValue 0x22 0x00 0xCC   # hlink
$cs.Colors(12).RGB = ColorValue 0x55 0x1A 0x8B   # folHlink
